$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.311628556340793
$ws.Range("C2").Value = 0.182196140785436
$ws.Range("D2").Value = 0.1501265540252845
$ws.Range("F2").Value = 1.699915982715396
$ws.Range("G2").Value = 0.002482491547813603
$ws.Range("J2").Value = 0.2035584618202839
$ws.Range("L2").Value = 0.3290868311613266
$ws.Range("M2").Value = 0.3294330702549146
$ws.Range("N2").Value = 1.692190539315551
$ws.Range("O2").Value = 4.306380649986352
$ws.Range("B3").Value = 1.233359189711052
$ws.Range("C3").Value = 0.172647344870839
$ws.Range("D3").Value = 0.150081574071816
$ws.Range("F3").Value = 1.704008030465303
$ws.Range("G3").Value = 0.002485416671872163
$ws.Range("J3").Value = 0.2049127701997739
$ws.Range("L3").Value = 0.326436158904464
$ws.Range("M3").Value = 0.31695269073753
$ws.Range("N3").Value = 1.708474304463831
$ws.Range("O3").Value = 4.317639526875553
$ws.Range("B4").Value = 1.185678597658182
$ws.Range("C4").Value = 0.166736482038587
$ws.Range("D4").Value = 0.1500908570871111
$ws.Range("F4").Value = 1.707379432134964
$ws.Range("G4").Value = 0.002487310398538159
$ws.Range("J4").Value = 0.2057890190080882
$ws.Range("L4").Value = 0.3249205585359292
$ws.Range("M4").Value = 0.3094023283990879
$ws.Range("N4").Value = 1.719058984731713
$ws.Range("O4").Value = 4.326877455815207
$ws.Range("B5").Value = 1.166344512356972
$ws.Range("C5").Value = 0.1643158433782617
$ws.Range("D5").Value = 0.1501039536213717
$ws.Range("F5").Value = 1.708969372386655
$ws.Range("G5").Value = 0.002488106745171897
$ws.Range("J5").Value = 0.206157354805826
$ws.Range("L5").Value = 0.3243311951768035
$ws.Range("M5").Value = 0.3063540625521668
$ws.Range("N5").Value = 1.723519853193842
$ws.Range("O5").Value = 4.331226640185349
$ws.Range("B6").Value = 1.163139947769025
$ws.Range("C6").Value = 0.1639131821803943
$ws.Range("D6").Value = 0.1501066918299045
$ws.Range("F6").Value = 1.709246433540343
$ws.Range("G6").Value = 0.0024882404682115
$ws.Range("J6").Value = 0.2062191971165737
$ws.Range("L6").Value = 0.3242350417417086
$ws.Range("M6").Value = 0.3058496329523166
$ws.Range("N6").Value = 1.724269489433226
$ws.Range("O6").Value = 4.331984131525871
$ws.Range("B7").Value = 1.18541746040043
$ws.Range("C7").Value = 0.1667038845381796
$ws.Range("D7").Value = 0.150090995954649
$ws.Range("F7").Value = 1.70739999965172
$ws.Range("G7").Value = 0.002487321038317783
$ws.Range("J7").Value = 0.2057939409094702
$ws.Range("L7").Value = 0.3249124956184488
$ws.Range("M7").Value = 0.3093611024255054
$ws.Range("N7").Value = 1.719118548139818
$ws.Range("O7").Value = 4.326933743163323
$ws.Range("B8").Value = 1.284563770202908
$ws.Range("C8").Value = 0.1789137202355988
$ws.Range("D8").Value = 0.1501034066652309
$ws.Range("F8").Value = 1.701148717829589
$ws.Range("G8").Value = 0.002483479899425965
$ws.Range("J8").Value = 0.2040161609837687
$ws.Range("L8").Value = 0.328149715495087
$ws.Range("M8").Value = 0.3251065982939849
$ws.Range("N8").Value = 1.697683523127761
$ws.Range("O8").Value = 4.309780179875418
$ws.Range("B9").Value = 1.48193431976631
$ws.Range("C9").Value = 0.2024731049094441
$ws.Range("D9").Value = 0.1504191080504143
$ws.Range("F9").Value = 1.695702025016232
$ws.Range("G9").Value = 0.002476719206930844
$ws.Range("J9").Value = 0.200883767061977
$ws.Range("L9").Value = 0.3353816144695116
$ws.Range("M9").Value = 0.3568683109327324
$ws.Range("N9").Value = 1.660298525510314
$ws.Range("O9").Value = 4.294592461885458
$ws.Range("B10").Value = 1.628689587167287
$ws.Range("C10").Value = 0.2195438328707269
$ws.Range("D10").Value = 0.1508268931251493
$ws.Range("F10").Value = 1.695851922873345
$ws.Range("G10").Value = 0.002472217934525076
$ws.Range("J10").Value = 0.1987968973753091
$ws.Range("L10").Value = 0.3412286630835268
$ws.Range("M10").Value = 0.3807337200312091
$ws.Range("N10").Value = 1.635660445497358
$ws.Range("O10").Value = 4.294692523385578
$ws.Range("B11").Value = 1.695822470510393
$ws.Range("C11").Value = 0.2272571938784154
$ws.Range("D11").Value = 0.151050219179325
$ws.Range("F11").Value = 1.696821512283208
$ws.Range("G11").Value = 0.002470270339603123
$ws.Range("J11").Value = 0.1978938462103061
$ws.Range("L11").Value = 0.3440035452459824
$ws.Range("M11").Value = 0.3917039986186595
$ws.Range("N11").Value = 1.625065057481152
$ws.Range("O11").Value = 4.297185473541276
$ws.Range("B12").Value = 1.721296488143992
$ws.Range("C12").Value = 0.2301704336480555
$ws.Range("D12").Value = 0.1511401929947169
$ws.Range("F12").Value = 1.697318244379986
$ws.Range("G12").Value = 0.002469547148666232
$ws.Range("J12").Value = 0.1975585205737662
$ws.Range("L12").Value = 0.3450707616437256
$ws.Range("M12").Value = 0.3958743096501038
$ws.Range("N12").Value = 1.621140895737149
$ws.Range("O12").Value = 4.298481546631592
$ws.Range("B13").Value = 1.715807901185769
$ws.Range("C13").Value = 0.2295433574189474
$ws.Range("D13").Value = 0.1511205756677896
$ws.Range("F13").Value = 1.6972055023965
$ws.Range("G13").Value = 0.002469702264873612
$ws.Range("J13").Value = 0.1976304438478014
$ws.Range("L13").Value = 0.3448401886128494
$ws.Range("M13").Value = 0.3949754463900987
$ws.Range("N13").Value = 1.621982116821428
$ws.Range("O13").Value = 4.298186754978275
$ws.Range("B14").Value = 1.697917194143201
$ws.Range("C14").Value = 0.2274970218195449
$ws.Range("D14").Value = 0.1510575132461796
$ws.Range("F14").Value = 1.696859782253185
$ws.Range("G14").Value = 0.002470210555626864
$ws.Range("J14").Value = 0.1978661257690362
$ws.Range("L14").Value = 0.3440910171678979
$ws.Range("M14").Value = 0.392046771329575
$ws.Range("N14").Value = 1.624740449080683
$ws.Range("O14").Value = 4.297285046083744
$ws.Range("B15").Value = 1.686965382380322
$ws.Range("C15").Value = 0.22624258283264
$ws.Range("D15").Value = 0.1510195885948633
$ws.Range("F15").Value = 1.696664890971775
$ws.Range("G15").Value = 0.002470523761401394
$ws.Range("J15").Value = 0.1980113520864757
$ws.Range("L15").Value = 0.3436342641999062
$ws.Range("M15").Value = 0.3902549629551544
$ws.Range("N15").Value = 1.62644147788005
$ws.Range("O15").Value = 4.296778573496908
$ws.Range("B16").Value = 1.624309683983654
$ws.Range("C16").Value = 0.2190386845387025
$ws.Range("D16").Value = 0.1508130564396808
$ws.Range("F16").Value = 1.695806692072964
$ws.Range("G16").Value = 0.002472347222108541
$ws.Range("J16").Value = 0.1988568438702398
$ws.Range("L16").Value = 0.3410496230918056
$ws.Range("M16").Value = 0.3800190541409094
$ws.Range("N16").Value = 1.63636520518579
$ws.Range("O16").Value = 4.294578861155543
$ws.Range("B17").Value = 1.585967017416351
$ws.Range("C17").Value = 0.2146058596867704
$ws.Range("D17").Value = 0.1506960188600814
$ws.Range("F17").Value = 1.695511049023438
$ws.Range("G17").Value = 0.002473491435261808
$ws.Range("J17").Value = 0.1993873679888569
$ws.Range("L17").Value = 0.3394934114115955
$ws.Range("M17").Value = 0.373768615205627
$ws.Range("N17").Value = 1.642609983697092
$ws.Range("O17").Value = 4.293856332645248
$ws.Range("B18").Value = 1.563948539005082
$ws.Range("C18").Value = 0.2120513110244531
$ws.Range("D18").Value = 0.1506322639907225
$ws.Range("F18").Value = 1.695425850022389
$ws.Range("G18").Value = 0.00247415897827242
$ws.Range("J18").Value = 0.1996968680176394
$ws.Range("L18").Value = 0.3386091535287932
$ws.Range("M18").Value = 0.370184253323103
$ws.Range("N18").Value = 1.64625947603048
$ws.Range("O18").Value = 4.293671110657499
$ws.Range("B19").Value = 1.556499550792068
$ws.Range("C19").Value = 0.2111855464583243
$ws.Range("D19").Value = 0.1506112905037966
$ws.Range("F19").Value = 1.695411577215538
$ws.Range("G19").Value = 0.002474386616873553
$ws.Range("J19").Value = 0.1998024080937375
$ws.Range("L19").Value = 0.3383116231707817
$ws.Range("M19").Value = 0.368972499905361
$ws.Range("N19").Value = 1.647505035094738
$ws.Range("O19").Value = 4.293647958444211
$ws.Range("B20").Value = 1.590045025923985
$ws.Range("C20").Value = 0.2150782501045967
$ws.Range("D20").Value = 0.1507081092967297
$ws.Range("F20").Value = 1.695533739742885
$ws.Range("G20").Value = 0.002473368657275644
$ws.Range("J20").Value = 0.1993304420037592
$ws.Range("L20").Value = 0.3396579523481051
$ws.Range("M20").Value = 0.3744328769410288
$ws.Range("N20").Value = 1.641939248753971
$ws.Range("O20").Value = 4.293909405123202
$ws.Range("B21").Value = 1.703170721613276
$ws.Range("C21").Value = 0.2280982886942695
$ws.Range("D21").Value = 0.1510758897836126
$ws.Range("F21").Value = 1.696957812560242
$ws.Range("G21").Value = 0.002470060870069108
$ws.Range("J21").Value = 0.1977967201744066
$ws.Range("L21").Value = 0.3443106220052385
$ws.Range("M21").Value = 0.3929065591617089
$ws.Range("N21").Value = 1.623927869571247
$ws.Range("O21").Value = 4.297540344519206
$ws.Range("B22").Value = 1.777408813112118
$ws.Range("C22").Value = 0.2365630262436866
$ws.Range("D22").Value = 0.1513477450000096
$ws.Range("F22").Value = 1.698643733972787
$ws.Range("G22").Value = 0.002467982483532429
$ws.Range("J22").Value = 0.1968330422100697
$ws.Range("L22").Value = 0.3474471178874836
$ws.Range("M22").Value = 0.405073929300336
$ws.Range("N22").Value = 1.612669798683427
$ws.Range("O22").Value = 4.301965423115291
$ws.Range("B23").Value = 1.737759204245208
$ws.Range("C23").Value = 0.2320493633926048
$ws.Range("D23").Value = 0.1511997805889749
$ws.Range("F23").Value = 1.697674841164513
$ws.Range("G23").Value = 0.002469084144618237
$ws.Range("J23").Value = 0.197343839077218
$ws.Range("L23").Value = 0.3457643908758001
$ws.Range("M23").Value = 0.3985714815971662
$ws.Range("N23").Value = 1.618631468168942
$ws.Range("O23").Value = 4.299415876626483
$ws.Range("B24").Value = 1.588201279575458
$ws.Range("C24").Value = 0.2148647012580795
$ws.Range("D24").Value = 0.1507026322089402
$ws.Range("F24").Value = 1.695523217209328
$ws.Range("G24").Value = 0.002473424135026603
$ws.Range("J24").Value = 0.1993561642402821
$ws.Range("L24").Value = 0.3395835308977695
$ws.Range("M24").Value = 0.3741325358492631
$ws.Range("N24").Value = 1.642242303345171
$ws.Range("O24").Value = 4.29388469409372
$ws.Range("B25").Value = 1.428230072109386
$ws.Range("C25").Value = 0.1961412208585784
$ws.Range("D25").Value = 0.1503026948927655
$ws.Range("F25").Value = 1.696446350676894
$ws.Range("G25").Value = 0.002478466016245298
$ws.Range("J25").Value = 0.201693407890299
$ws.Range("L25").Value = 0.3353816144695116
$ws.Range("M25").Value = 0.3568683109327324
$ws.Range("N25").Value = 1.669915032798571
$ws.Range("O25").Value = 4.296724806145562
